$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1952861952861953
$ws.Range("C2").Value = 0.569023569023569
$ws.Range("J2").Value = 0.01683501683501683
$ws.Range("P2").Value = 0.1582491582491583
$ws.Range("S2").Value = 0.06060606060606061
$ws.Range("B3").Value = 0.02162162162162162
$ws.Range("C3").Value = 0.03783783783783784
$ws.Range("J3").Value = 0.05945945945945946
$ws.Range("P3").Value = 0.7351351351351352
$ws.Range("S3").Value = 0.145945945945946
$ws.Range("J4").Value = 0.1025641025641026
$ws.Range("P4").Value = 0.6153846153846154
$ws.Range("S4").Value = 0.282051282051282
$ws.Range("B6").Value = 0.07017543859649122
$ws.Range("D6").Value = 0.008771929824561403
$ws.Range("F6").Value = 0.05263157894736842
$ws.Range("J6").Value = 0.2763157894736842
$ws.Range("O6").Value = 0.0131578947368421
$ws.Range("Q6").Value = 0.1140350877192982
$ws.Range("R6").Value = 0.08333333333333333
$ws.Range("S6").Value = 0.3815789473684211
$ws.Range("B7").Value = 0.1421319796954315
$ws.Range("D7").Value = 0.01522842639593909
$ws.Range("F7").Value = 0.05076142131979695
$ws.Range("J7").Value = 0.1421319796954315
$ws.Range("O7").Value = 0.005076142131979695
$ws.Range("Q7").Value = 0.1624365482233502
$ws.Range("R7").Value = 0.06598984771573604
$ws.Range("S7").Value = 0.4162436548223351
$ws.Range("B8").Value = 0.1191011235955056
$ws.Range("D8").Value = 0.01348314606741573
$ws.Range("E8").Value = 0.002247191011235955
$ws.Range("F8").Value = 0.0449438202247191
$ws.Range("J8").Value = 0.1078651685393258
$ws.Range("O8").Value = 0.01123595505617977
$ws.Range("Q8").Value = 0.2067415730337079
$ws.Range("R8").Value = 0.07865168539325842
$ws.Range("S8").Value = 0.4157303370786517
$ws.Range("B9").Value = 0.119047619047619
$ws.Range("D9").Value = 0.01587301587301587
$ws.Range("E9").Value = 0.007936507936507936
$ws.Range("F9").Value = 0.07936507936507936
$ws.Range("Q9").Value = 0.1507936507936508
$ws.Range("R9").Value = 0.07936507936507936
$ws.Range("S9").Value = 0.4047619047619048
$ws.Range("B10").Value = 0.08969315499606609
$ws.Range("D10").Value = 0.02281667977970102
$ws.Range("E10").Value = 0.0007867820613690008
$ws.Range("F10").Value = 0.07238394964594808
$ws.Range("J10").Value = 0.1337529504327301
$ws.Range("O10").Value = 0.01494885916601102
$ws.Range("Q10").Value = 0.2289535798583792
$ws.Range("R10").Value = 0.07395751376868608
$ws.Range("S10").Value = 0.3627065302911094
$ws.Range("G11").Value = 0.1485148514851485
$ws.Range("J11").Value = 0.0924092409240924
$ws.Range("K11").Value = 0.2211221122112211
$ws.Range("L11").Value = 0.5313531353135313
$ws.Range("S11").Value = 0.006600660066006601
$ws.Range("G12").Value = 0.7530120481927711
$ws.Range("J12").Value = 0.1927710843373494
$ws.Range("K12").Value = 0.006024096385542169
$ws.Range("L12").Value = 0.03012048192771084
$ws.Range("S12").Value = 0.01807228915662651
$ws.Range("G13").Value = 0.6938775510204082
$ws.Range("J13").Value = 0.3061224489795918
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 0.02941176470588235
$ws.Range("H15").Value = 0.142156862745098
$ws.Range("I15").Value = 0.04411764705882353
$ws.Range("J15").Value = 0.4019607843137255
$ws.Range("K15").Value = 0.107843137254902
$ws.Range("M15").Value = 0.01470588235294118
$ws.Range("N15").Value = 0.004901960784313725
$ws.Range("O15").Value = 0.07352941176470588
$ws.Range("S15").Value = 0.1813725490196078
$ws.Range("F16").Value = 0.03465346534653466
$ws.Range("H16").Value = 0.1485148514851485
$ws.Range("I16").Value = 0.06930693069306931
$ws.Range("J16").Value = 0.4405940594059406
$ws.Range("K16").Value = 0.09900990099009901
$ws.Range("M16").Value = 0.0198019801980198
$ws.Range("O16").Value = 0.06930693069306931
$ws.Range("S16").Value = 0.1188118811881188
$ws.Range("F17").Value = 0.03036876355748373
$ws.Range("H17").Value = 0.1626898047722343
$ws.Range("I17").Value = 0.06941431670281996
$ws.Range("J17").Value = 0.4425162689804772
$ws.Range("K17").Value = 0.1106290672451193
$ws.Range("M17").Value = 0.01952277657266811
$ws.Range("O17").Value = 0.04772234273318872
$ws.Range("S17").Value = 0.1171366594360087
$ws.Range("F18").Value = 0.01142857142857143
$ws.Range("I18").Value = 0.02857142857142857
$ws.Range("J18").Value = 0.4342857142857143
$ws.Range("K18").Value = 0.05714285714285714
$ws.Range("M18").Value = 0.02857142857142857
$ws.Range("O18").Value = 0.08
$ws.Range("S18").Value = 0.16
$ws.Range("F19").Value = 0.02113271344040575
$ws.Range("H19").Value = 0.231614539306847
$ws.Range("J19").Value = 0.3609467455621302
$ws.Range("K19").Value = 0.1149619611158073
$ws.Range("M19").Value = 0.0253592561284869
$ws.Range("O19").Value = 0.073541842772612
$ws.Range("S19").Value = 0.117497886728656
